$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 36,16
$data[0,0] = 3
$data[0,1] = 1
$data[0,2] = 8.723857666666666
$data[0,3] = 26.171573
$data[0,4] = 0.007098432040951201
$data[0,5] = 0.007098432040951203
$data[0,6] = 3
$data[0,7] = 1
$data[0,8] = 2.092292333333333
$data[0,9] = 6.276877000000001
$data[0,10] = 0.1237967521619938
$data[0,11] = 0.1237967521619938
$data[0,12] = 18.25286051305789
$data[0,13] = 164.275744617521
$data[0,14] = 0.0008787628321123914
$data[0,15] = 0.0008787628321123917
$data[1,0] = 3
$data[1,1] = 1
$data[1,2] = 8.723857666666666
$data[1,3] = 26.171573
$data[1,4] = 0.007098432040951201
$data[1,5] = 0.007098432040951203
$data[1,6] = 3
$data[1,7] = 1
$data[1,8] = 2.468365333333333
$data[1,9] = 7.405096
$data[1,10] = 0.1460482393151517
$data[1,11] = 0.1460482393151517
$data[1,12] = 21.53366783733422
$data[1,13] = 193.803010536008
$data[1,14] = 0.001036713501479182
$data[1,15] = 0.001036713501479182
$data[2,0] = 3
$data[2,1] = 1
$data[2,2] = 8.723857666666666
$data[2,3] = 26.171573
$data[2,4] = 0.007098432040951201
$data[2,5] = 0.007098432040951203
$data[2,6] = 3
$data[2,7] = 1
$data[2,8] = 0.315935
$data[2,9] = 0.9478049999999999
$data[2,10] = 0.01869324198688273
$data[2,11] = 0.01869324198688273
$data[2,12] = 2.756171971918333
$data[2,13] = 24.805547747265
$data[2,14] = 0.0001326927078689426
$data[2,15] = 0.0001326927078689427
$data[3,0] = 3
$data[3,1] = 1
$data[3,2] = 8.723857666666666
$data[3,3] = 26.171573
$data[3,4] = 0.007098432040951201
$data[3,5] = 0.007098432040951203
$data[3,6] = 3
$data[3,7] = 1
$data[3,8] = 0.4705663333333334
$data[3,9] = 1.411699
$data[3,10] = 0.02784246867197405
$data[3,11] = 0.02784246867197405
$data[3,12] = 4.105153714725222
$data[3,13] = 36.946383432527
$data[3,14] = 0.0001976378717203206
$data[3,15] = 0.0001976378717203207
$data[4,0] = 3
$data[4,1] = 1
$data[4,2] = 8.723857666666666
$data[4,3] = 26.171573
$data[4,4] = 0.007098432040951201
$data[4,5] = 0.007098432040951203
$data[4,6] = 3
$data[4,7] = 1
$data[4,8] = 11.291786
$data[4,9] = 33.875358
$data[4,10] = 0.6681123907199095
$data[4,11] = 0.6681123907199095
$data[4,12] = 98.50793386645933
$data[4,13] = 886.5714047981339
$data[4,14] = 0.004742550401242713
$data[4,15] = 0.004742550401242714
$data[5,0] = 3
$data[5,1] = 1
$data[5,2] = 8.723857666666666
$data[5,3] = 26.171573
$data[5,4] = 0.007098432040951201
$data[5,5] = 0.007098432040951203
$data[5,6] = 3
$data[5,7] = 1
$data[5,8] = 0.2620826666666667
$data[5,9] = 0.7862480000000001
$data[5,10] = 0.01550690714408826
$data[5,11] = 0.01550690714408826
$data[5,12] = 2.286371880900445
$data[5,13] = 20.577346928104
$data[5,14] = 0.0001100747265276512
$data[5,15] = 0.0001100747265276512
$data[6,0] = 3
$data[6,1] = 1
$data[6,2] = 1065.000325333333
$data[6,3] = 3195.000976
$data[6,4] = 0.8665698962346957
$data[6,5] = 0.8665698962346958
$data[6,6] = 3
$data[6,7] = 1
$data[6,8] = 2.092292333333333
$data[6,9] = 6.276877000000001
$data[6,10] = 0.1237967521619938
$data[6,11] = 0.1237967521619938
$data[6,12] = 2228.292015692439
$data[6,13] = 20054.62814123195
$data[6,14] = 0.1072785386752113
$data[6,15] = 0.1072785386752113
$data[7,0] = 3
$data[7,1] = 1
$data[7,2] = 1065.000325333333
$data[7,3] = 3195.000976
$data[7,4] = 0.8665698962346957
$data[7,5] = 0.8665698962346958
$data[7,6] = 3
$data[7,7] = 1
$data[7,8] = 2.468365333333333
$data[7,9] = 7.405096
$data[7,10] = 0.1460482393151517
$data[7,11] = 0.1460482393151517
$data[7,12] = 2628.809883041522
$data[7,13] = 23659.2889473737
$data[7,14] = 0.126561007588591
$data[7,15] = 0.126561007588591
$data[8,0] = 3
$data[8,1] = 1
$data[8,2] = 1065.000325333333
$data[8,3] = 3195.000976
$data[8,4] = 0.8665698962346957
$data[8,5] = 0.8665698962346958
$data[8,6] = 3
$data[8,7] = 1
$data[8,8] = 0.315935
$data[8,9] = 0.9478049999999999
$data[8,10] = 0.01869324198688273
$data[8,11] = 0.01869324198688273
$data[8,12] = 336.4708777841867
$data[8,13] = 3028.23790005768
$data[8,14] = 0.01619900076886302
$data[8,15] = 0.01619900076886302
$data[9,0] = 3
$data[9,1] = 1
$data[9,2] = 1065.000325333333
$data[9,3] = 3195.000976
$data[9,4] = 0.8665698962346957
$data[9,5] = 0.8665698962346958
$data[9,6] = 3
$data[9,7] = 1
$data[9,8] = 0.4705663333333334
$data[9,9] = 1.411699
$data[9,10] = 0.02784246867197405
$data[9,11] = 0.02784246867197405
$data[9,12] = 501.1532980909138
$data[9,13] = 4510.379682818225
$data[9,14] = 0.02412744518799032
$data[9,15] = 0.02412744518799032
$data[10,0] = 3
$data[10,1] = 1
$data[10,2] = 1065.000325333333
$data[10,3] = 3195.000976
$data[10,4] = 0.8665698962346957
$data[10,5] = 0.8665698962346958
$data[10,6] = 3
$data[10,7] = 1
$data[10,8] = 11.291786
$data[10,9] = 33.875358
$data[10,10] = 0.6681123907199095
$data[10,11] = 0.6681123907199095
$data[10,12] = 12025.75576359438
$data[10,13] = 108231.8018723494
$data[10,14] = 0.5789660850992664
$data[10,15] = 0.5789660850992665
$data[11,0] = 3
$data[11,1] = 1
$data[11,2] = 1065.000325333333
$data[11,3] = 3195.000976
$data[11,4] = 0.8665698962346957
$data[11,5] = 0.8665698962346958
$data[11,6] = 3
$data[11,7] = 1
$data[11,8] = 0.2620826666666667
$data[11,9] = 0.7862480000000001
$data[11,10] = 0.01550690714408826
$data[11,11] = 0.01550690714408826
$data[11,12] = 279.1181252642276
$data[11,13] = 2512.063127378048
$data[11,14] = 0.01343781891477363
$data[11,15] = 0.01343781891477363
$data[12,0] = 3
$data[12,1] = 1
$data[12,2] = 0.8450703333333333
$data[12,3] = 2.535211
$data[12,4] = 0.0006876171712327699
$data[12,5] = 0.0006876171712327699
$data[12,6] = 3
$data[12,7] = 1
$data[12,8] = 2.092292333333333
$data[12,9] = 6.276877000000001
$data[12,10] = 0.1237967521619938
$data[12,11] = 0.1237967521619938
$data[12,12] = 1.768134179560778
$data[12,13] = 15.913207616047
$data[12,14] = 0.00008512477252943445
$data[12,15] = 0.00008512477252943446
$data[13,0] = 3
$data[13,1] = 1
$data[13,2] = 0.8450703333333333
$data[13,3] = 2.535211
$data[13,4] = 0.0006876171712327699
$data[13,5] = 0.0006876171712327699
$data[13,6] = 3
$data[13,7] = 1
$data[13,8] = 2.468365333333333
$data[13,9] = 7.405096
$data[13,10] = 0.1460482393151517
$data[13,11] = 0.1460482393151517
$data[13,12] = 2.085942315028444
$data[13,13] = 18.773480835256
$data[13,14] = 0.0001004252771814112
$data[13,15] = 0.0001004252771814112
$data[14,0] = 3
$data[14,1] = 1
$data[14,2] = 0.8450703333333333
$data[14,3] = 2.535211
$data[14,4] = 0.0006876171712327699
$data[14,5] = 0.0006876171712327699
$data[14,6] = 3
$data[14,7] = 1
$data[14,8] = 0.315935
$data[14,9] = 0.9478049999999999
$data[14,10] = 0.01869324198688273
$data[14,11] = 0.01869324198688273
$data[14,12] = 0.2669872957616666
$data[14,13] = 2.402885661855
$data[14,14] = 0.00001285379417618994
$data[14,15] = 0.00001285379417618994
$data[15,0] = 3
$data[15,1] = 1
$data[15,2] = 0.8450703333333333
$data[15,3] = 2.535211
$data[15,4] = 0.0006876171712327699
$data[15,5] = 0.0006876171712327699
$data[15,6] = 3
$data[15,7] = 1
$data[15,8] = 0.4705663333333334
$data[15,9] = 1.411699
$data[15,10] = 0.02784246867197405
$data[15,11] = 0.02784246867197405
$data[15,12] = 0.3976616481654444
$data[15,13] = 3.578954833489
$data[15,14] = 0.00001914495954835981
$data[15,15] = 0.00001914495954835981
$data[16,0] = 3
$data[16,1] = 1
$data[16,2] = 0.8450703333333333
$data[16,3] = 2.535211
$data[16,4] = 0.0006876171712327699
$data[16,5] = 0.0006876171712327699
$data[16,6] = 3
$data[16,7] = 1
$data[16,8] = 11.291786
$data[16,9] = 33.875358
$data[16,10] = 0.6681123907199095
$data[16,11] = 0.6681123907199095
$data[16,12] = 9.542353358948667
$data[16,13] = 85.88118023053799
$data[16,14] = 0.0004594055521723873
$data[16,15] = 0.0004594055521723873
$data[17,0] = 3
$data[17,1] = 1
$data[17,2] = 0.8450703333333333
$data[17,3] = 2.535211
$data[17,4] = 0.0006876171712327699
$data[17,5] = 0.0006876171712327699
$data[17,6] = 3
$data[17,7] = 1
$data[17,8] = 0.2620826666666667
$data[17,9] = 0.7862480000000001
$data[17,10] = 0.01550690714408826
$data[17,11] = 0.01550690714408826
$data[17,12] = 0.2214782864808889
$data[17,13] = 1.993304578328
$data[17,14] = 0.0000106628156249872
$data[17,15] = 0.0000106628156249872
$data[18,0] = 3
$data[18,1] = 1
$data[18,2] = 1.780731666666667
$data[18,3] = 5.342195
$data[18,4] = 0.001448946464051256
$data[18,5] = 0.001448946464051256
$data[18,6] = 3
$data[18,7] = 1
$data[18,8] = 2.092292333333333
$data[18,9] = 6.276877000000001
$data[18,10] = 0.1237967521619938
$data[18,11] = 0.1237967521619938
$data[18,12] = 3.725811213890556
$data[18,13] = 33.53230092501501
$data[18,14] = 0.0001793748663061505
$data[18,15] = 0.0001793748663061505
$data[19,0] = 3
$data[19,1] = 1
$data[19,2] = 1.780731666666667
$data[19,3] = 5.342195
$data[19,4] = 0.001448946464051256
$data[19,5] = 0.001448946464051256
$data[19,6] = 3
$data[19,7] = 1
$data[19,8] = 2.468365333333333
$data[19,9] = 7.405096
$data[19,10] = 0.1460482393151517
$data[19,11] = 0.1460482393151517
$data[19,12] = 4.39549631396889
$data[19,13] = 39.55946682572
$data[19,14] = 0.0002116160799366006
$data[19,15] = 0.0002116160799366006
$data[20,0] = 3
$data[20,1] = 1
$data[20,2] = 1.780731666666667
$data[20,3] = 5.342195
$data[20,4] = 0.001448946464051256
$data[20,5] = 0.001448946464051256
$data[20,6] = 3
$data[20,7] = 1
$data[20,8] = 0.315935
$data[20,9] = 0.9478049999999999
$data[20,10] = 0.01869324198688273
$data[20,11] = 0.01869324198688273
$data[20,12] = 0.5625954591083333
$data[20,13] = 5.063359131975
$data[20,14] = 0.00002708550687854819
$data[20,15] = 0.00002708550687854819
$data[21,0] = 3
$data[21,1] = 1
$data[21,2] = 1.780731666666667
$data[21,3] = 5.342195
$data[21,4] = 0.001448946464051256
$data[21,5] = 0.001448946464051256
$data[21,6] = 3
$data[21,7] = 1
$data[21,8] = 0.4705663333333334
$data[21,9] = 1.411699
$data[21,10] = 0.02784246867197405
$data[21,11] = 0.02784246867197405
$data[21,12] = 0.837952371033889
$data[21,13] = 7.541571339305
$data[21,14] = 0.00004034224653271465
$data[21,15] = 0.00004034224653271466
$data[22,0] = 3
$data[22,1] = 1
$data[22,2] = 1.780731666666667
$data[22,3] = 5.342195
$data[22,4] = 0.001448946464051256
$data[22,5] = 0.001448946464051256
$data[22,6] = 3
$data[22,7] = 1
$data[22,8] = 11.291786
$data[22,9] = 33.875358
$data[22,10] = 0.6681123907199095
$data[22,11] = 0.6681123907199095
$data[22,12] = 20.10764090342333
$data[22,13] = 180.96876813081
$data[22,14] = 0.0009680590861224437
$data[22,15] = 0.0009680590861224437
$data[23,0] = 3
$data[23,1] = 1
$data[23,2] = 1.780731666666667
$data[23,3] = 5.342195
$data[23,4] = 0.001448946464051256
$data[23,5] = 0.001448946464051256
$data[23,6] = 3
$data[23,7] = 1
$data[23,8] = 0.2620826666666667
$data[23,9] = 0.7862480000000001
$data[23,10] = 0.01550690714408826
$data[23,11] = 0.01550690714408826
$data[23,12] = 0.4666989038177778
$data[23,13] = 4.20029013436
$data[23,14] = 0.00002246867827479784
$data[23,15] = 0.00002246867827479784
$data[24,0] = 3
$data[24,1] = 1
$data[24,2] = 0.2508443333333333
$data[24,3] = 0.752533
$data[24,4] = 0.0002041071187839237
$data[24,5] = 0.0002041071187839237
$data[24,6] = 3
$data[24,7] = 1
$data[24,8] = 2.092292333333333
$data[24,9] = 6.276877000000001
$data[24,10] = 0.1237967521619938
$data[24,11] = 0.1237967521619938
$data[24,12] = 0.5248396754934445
$data[24,13] = 4.723557079441001
$data[24,14] = 0.00002526779839859203
$data[24,15] = 0.00002526779839859203
$data[25,0] = 3
$data[25,1] = 1
$data[25,2] = 0.2508443333333333
$data[25,3] = 0.752533
$data[25,4] = 0.0002041071187839237
$data[25,5] = 0.0002041071187839237
$data[25,6] = 3
$data[25,7] = 1
$data[25,8] = 2.468365333333333
$data[25,9] = 7.405096
$data[25,10] = 0.1460482393151517
$data[25,11] = 0.1460482393151517
$data[25,12] = 0.6191754564631111
$data[25,13] = 5.572579108168
$data[25,14] = 0.00002980948533008059
$data[25,15] = 0.00002980948533008059
$data[26,0] = 3
$data[26,1] = 1
$data[26,2] = 0.2508443333333333
$data[26,3] = 0.752533
$data[26,4] = 0.0002041071187839237
$data[26,5] = 0.0002041071187839237
$data[26,6] = 3
$data[26,7] = 1
$data[26,8] = 0.315935
$data[26,9] = 0.9478049999999999
$data[26,10] = 0.01869324198688273
$data[26,11] = 0.01869324198688273
$data[26,12] = 0.07925050445166666
$data[26,13] = 0.7132545400649999
$data[26,14] = 0.0000038154237626733
$data[26,15] = 0.0000038154237626733
$data[27,0] = 3
$data[27,1] = 1
$data[27,2] = 0.2508443333333333
$data[27,3] = 0.752533
$data[27,4] = 0.0002041071187839237
$data[27,5] = 0.0002041071187839237
$data[27,6] = 3
$data[27,7] = 1
$data[27,8] = 0.4705663333333334
$data[27,9] = 1.411699
$data[27,10] = 0.02784246867197405
$data[27,11] = 0.02784246867197405
$data[27,12] = 0.1180388981741111
$data[27,13] = 1.062350083567
$data[27,14] = 0.00000568284606046828
$data[27,15] = 0.00000568284606046828
$data[28,0] = 3
$data[28,1] = 1
$data[28,2] = 0.2508443333333333
$data[28,3] = 0.752533
$data[28,4] = 0.0002041071187839237
$data[28,5] = 0.0002041071187839237
$data[28,6] = 3
$data[28,7] = 1
$data[28,8] = 11.291786
$data[28,9] = 33.875358
$data[28,10] = 0.6681123907199095
$data[28,11] = 0.6681123907199095
$data[28,12] = 2.832480531312667
$data[28,13] = 25.492324781814
$data[28,14] = 0.0001363664950936798
$data[28,15] = 0.0001363664950936798
$data[29,0] = 3
$data[29,1] = 1
$data[29,2] = 0.2508443333333333
$data[29,3] = 0.752533
$data[29,4] = 0.0002041071187839237
$data[29,5] = 0.0002041071187839237
$data[29,6] = 3
$data[29,7] = 1
$data[29,8] = 0.2620826666666667
$data[29,9] = 0.7862480000000001
$data[29,10] = 0.01550690714408826
$data[29,11] = 0.01550690714408826
$data[29,12] = 0.06574195179822223
$data[29,13] = 0.591677566184
$data[29,14] = 0.0000031650701384297
$data[29,15] = 0.0000031650701384297
$data[30,0] = 3
$data[30,1] = 1
$data[30,2] = 152.382926
$data[30,3] = 457.148778
$data[30,4] = 0.1239910009702851
$data[30,5] = 0.1239910009702851
$data[30,6] = 3
$data[30,7] = 1
$data[30,8] = 2.092292333333333
$data[30,9] = 6.276877000000001
$data[30,10] = 0.1237967521619938
$data[30,11] = 0.1237967521619938
$data[30,12] = 318.8296278007007
$data[30,13] = 2869.466650206306
$data[30,14] = 0.01534968321743592
$data[30,15] = 0.01534968321743592
$data[31,0] = 3
$data[31,1] = 1
$data[31,2] = 152.382926
$data[31,3] = 457.148778
$data[31,4] = 0.1239910009702851
$data[31,5] = 0.1239910009702851
$data[31,6] = 3
$data[31,7] = 1
$data[31,8] = 2.468365333333333
$data[31,9] = 7.405096
$data[31,10] = 0.1460482393151517
$data[31,11] = 0.1460482393151517
$data[31,12] = 376.1367319302987
$data[31,13] = 3385.230587372688
$data[31,14] = 0.01810866738263341
$data[31,15] = 0.01810866738263341
$data[32,0] = 3
$data[32,1] = 1
$data[32,2] = 152.382926
$data[32,3] = 457.148778
$data[32,4] = 0.1239910009702851
$data[32,5] = 0.1239910009702851
$data[32,6] = 3
$data[32,7] = 1
$data[32,8] = 0.315935
$data[32,9] = 0.9478049999999999
$data[32,10] = 0.01869324198688273
$data[32,11] = 0.01869324198688273
$data[32,12] = 48.14309972581
$data[32,13] = 433.2878975322899
$data[32,14] = 0.002317793785333351
$data[32,15] = 0.002317793785333351
$data[33,0] = 3
$data[33,1] = 1
$data[33,2] = 152.382926
$data[33,3] = 457.148778
$data[33,4] = 0.1239910009702851
$data[33,5] = 0.1239910009702851
$data[33,6] = 3
$data[33,7] = 1
$data[33,8] = 0.4705663333333334
$data[33,9] = 1.411699
$data[33,10] = 0.02784246867197405
$data[33,11] = 0.02784246867197405
$data[33,12] = 71.70627475042467
$data[33,13] = 645.356472753822
$data[33,14] = 0.003452215560121867
$data[33,15] = 0.003452215560121868
$data[34,0] = 3
$data[34,1] = 1
$data[34,2] = 152.382926
$data[34,3] = 457.148778
$data[34,4] = 0.1239910009702851
$data[34,5] = 0.1239910009702851
$data[34,6] = 3
$data[34,7] = 1
$data[34,8] = 11.291786
$data[34,9] = 33.875358
$data[34,10] = 0.6681123907199095
$data[34,11] = 0.6681123907199095
$data[34,12] = 1720.675390445836
$data[34,13] = 15486.07851401252
$data[34,14] = 0.0828399240860118
$data[34,15] = 0.08283992408601182
$data[35,0] = 3
$data[35,1] = 1
$data[35,2] = 152.382926
$data[35,3] = 457.148778
$data[35,4] = 0.1239910009702851
$data[35,5] = 0.1239910009702851
$data[35,6] = 3
$data[35,7] = 1
$data[35,8] = 0.2620826666666667
$data[35,9] = 0.7862480000000001
$data[35,10] = 0.01550690714408826
$data[35,11] = 0.01550690714408826
$data[35,12] = 39.93692360054933
$data[35,13] = 359.432312404944
$data[35,14] = 0.001922716938748769
$data[35,15] = 0.001922716938748769

$ws.Range("E2:T37").Value = $data
Write-Output "Done setting values"